$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (15.42578125 -> 16.42578125 in the saved OOXML).
# Excel's ColumnWidth COM property is quantized to a whole-pixel grid, so
# the nearest settable value that serializes closest to 16.42578125 is used.
$ws.Columns.Item(1).ColumnWidth = 15.65

# Update values in A1:A33
$ws.Cells.Item(1, 1).Value = 0.07640048516454101
$ws.Cells.Item(2, 1).Value = 0.016559779385612927
$ws.Cells.Item(3, 1).Value = -0.0039999999511142192
$ws.Cells.Item(4, 1).Value = -0.0079999999072501993
$ws.Cells.Item(5, 1).Value = -0.0029999999563505853
$ws.Cells.Item(6, 1).Value = -0.001999999960380805
$ws.Cells.Item(7, 1).Value = -0.0099999998771291843
$ws.Cells.Item(8, 1).Value = -0.0099999998765811782
$ws.Cells.Item(9, 1).Value = -0.0019999999602662299
$ws.Cells.Item(10, 1).Value = -0.0019999999618676156
$ws.Cells.Item(11, 1).Value = -0.0029999999518119935
$ws.Cells.Item(12, 1).Value = -0.0034999999473628307
$ws.Cells.Item(13, 1).Value = -0.0034999999523490644
$ws.Cells.Item(14, 1).Value = -0.007999999907935873
$ws.Cells.Item(15, 1).Value = -0.00099999998329636242
$ws.Cells.Item(16, 1).Value = 0.035900829050090177
$ws.Cells.Item(17, 1).Value = -0.0019999999755446751
$ws.Cells.Item(18, 1).Value = -0.0039999999544866327
$ws.Cells.Item(19, 1).Value = -0.0039999999579412027
$ws.Cells.Item(20, 1).Value = -0.0039999999542210674
$ws.Cells.Item(21, 1).Value = -0.0039999999537077002
$ws.Cells.Item(22, 1).Value = -0.0039999999533018027
$ws.Cells.Item(23, 1).Value = 0.0087829678395010546
$ws.Cells.Item(24, 1).Value = -0.019999999775739852
$ws.Cells.Item(25, 1).Value = -0.019999999772267962
$ws.Cells.Item(26, 1).Value = -0.002499999954761023
$ws.Cells.Item(27, 1).Value = -0.00249999995422856
$ws.Cells.Item(28, 1).Value = 0.059818684529216348
$ws.Cells.Item(29, 1).Value = -0.0069999998936767227
$ws.Cells.Item(30, 1).Value = -0.0599999993409841
$ws.Cells.Item(31, 1).Value = -0.0069999998879364256
$ws.Cells.Item(32, 1).Value = 0.012366786453140932
$ws.Cells.Item(33, 1).Value = -0.0039999999191113744
